$d = $word.ActiveDocument

# --- locate the insertion point -------------------------------------------------
# The new content belongs right after the "UVa 12405 - Scarecrow" / "Accepted"
# pair (both are interval-covering problems), and right before the empty
# paragraph that follows it.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "UVa 12405") {
        $anchor = $d.Paragraphs.Item($i + 1)
        break
    }
}

$pkgHeader = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# --- paragraph 1: "UVa 12321 - Gas Station (interval covering) " + " +" ---------
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null
$p1 = $d.Paragraphs.Item($anchor.Range.ListFormat.Count + 0)
$p1 = $anchor.Next()
$xml1 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-AU'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-AU'/></w:rPr><w:t xml:space='preserve'>UVa 12321 - Gas Station (interval covering) </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii='MS Gothic' w:eastAsia='MS Gothic' w:hAnsi='MS Gothic' w:cs='MS Gothic' w:hint='eastAsia'/><w:lang w:val='en-AU'/></w:rPr><w:t xml:space='preserve'> +</w:t></w:r></w:p>" + $pkgFooter
$p1.Range.InsertXML($xml1) | Out-Null

# --- paragraph 2: "Failing some uDebug cases" (ilvl 1) --------------------------
$r2 = $p1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter() | Out-Null
$p2 = $p1.Next()
$xml2 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-AU'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-AU'/></w:rPr><w:t>Failing some uDebug cases</w:t></w:r></w:p>" + $pkgFooter
$p2.Range.InsertXML($xml2) | Out-Null

# --- paragraph 3: empty paragraph that will host the relocated bookmark --------
$r3 = $p2.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter() | Out-Null
$p3 = $p2.Next()
$xml3 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-AU'/></w:rPr></w:pPr></w:p>" + $pkgFooter
$p3.Range.InsertXML($xml3) | Out-Null

# --- move the _GoBack bookmark from the last paragraph to the new paragraph ----
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$p3.Range.Bookmarks.Add("_GoBack") | Out-Null
